$wb = $excel.ActiveWorkbook

$signupSheet = $wb.Worksheets.Item("Signup")

# Add a new worksheet for Login data, placed after Signup
$loginSheet = $wb.Worksheets.Add($null, $signupSheet)
$loginSheet.Name = "Login"

# Column widths to match the Signup sheet's style
$loginSheet.Columns.Item(1).ColumnWidth = 14.67
$loginSheet.Columns.Item(2).ColumnWidth = 12.83

# Header row
$loginSheet.Range("A1").Value = "Email"
$loginSheet.Range("B1").Value = "Password"

# Data row (login credentials used for data-driven testing)
$loginSheet.Range("A2").Value = "test@gmail.com"
$loginSheet.Range("B2").Value = "Test@12345"

# Hyperlinks + hyperlink styling on the data row, mirroring the Signup sheet
$loginSheet.Hyperlinks.Add($loginSheet.Range("A2"), "mailto:test@gmail.com")
$loginSheet.Hyperlinks.Add($loginSheet.Range("B2"), "mailto:Test@12345")

# Restore the Signup sheet's selection
$signupSheet.Range("A5").Select() | Out-Null

# Make Login the active (selected) tab and set its selection
$loginSheet.Activate() | Out-Null
$loginSheet.Range("B11").Select() | Out-Null
